$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.829945333333333
$ws.Range("H2").Value = 5.489835999999999
$ws.Range("I2").Value = 0.4190796720210465
$ws.Range("J2").Value = 0.4190796720210465
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06447966666666667
$ws.Range("N2").Value = 0.193439
$ws.Range("O2").Value = 0.001101138907643723
$ws.Range("P2").Value = 0.001101138907643722
$ws.Range("Q2").Value = 0.1179942651115555
$ws.Range("R2").Value = 1.061948386004
$ws.Range("S2").Value = 0.0004614649322649446
$ws.Range("T2").Value = 0.0004614649322649445
$ws.Range("G3").Value = 1.829945333333333
$ws.Range("H3").Value = 5.489835999999999
$ws.Range("I3").Value = 0.4190796720210465
$ws.Range("J3").Value = 0.4190796720210465
$ws.Range("O3").Value = 0.00657695954769643
$ws.Range("P3").Value = 0.006576959547696431
$ws.Range("Q3").Value = 0.7047644062995555
$ws.Range("R3").Value = 6.342879656696
$ws.Range("S3").Value = 0.00275627005014431
$ws.Range("T3").Value = 0.00275627005014431
$ws.Range("G4").Value = 1.829945333333333
$ws.Range("H4").Value = 5.489835999999999
$ws.Range("I4").Value = 0.4190796720210465
$ws.Range("J4").Value = 0.4190796720210465
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008175
$ws.Range("N4").Value = 0.024525
$ws.Range("O4").Value = 0.0001396069650378791
$ws.Range("P4").Value = 0.0001396069650378791
$ws.Range("Q4").Value = 0.0149598031
$ws.Range("R4").Value = 0.1346382279
$ws.Range("S4").Value = 0.00005850644111992807
$ws.Range("T4").Value = 0.00005850644111992807
$ws.Range("G5").Value = 1.829945333333333
$ws.Range("H5").Value = 5.489835999999999
$ws.Range("I5").Value = 0.4190796720210465
$ws.Range("J5").Value = 0.4190796720210465
$ws.Range("M5").Value = 58.099467
$ws.Range("N5").Value = 174.298401
$ws.Range("O5").Value = 0.992182294579622
$ws.Range("P5").Value = 0.992182294579622
$ws.Range("Q5").Value = 106.318848505804
$ws.Range("R5").Value = 956.869636552236
$ws.Range("S5").Value = 0.4158034305975173
$ws.Range("T5").Value = 0.4158034305975173
$ws.Range("I6").Value = 0.2833335737960661
$ws.Range("J6").Value = 0.2833335737960661
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06447966666666667
$ws.Range("N6").Value = 0.193439
$ws.Range("O6").Value = 0.001101138907643723
$ws.Range("P6").Value = 0.001101138907643722
$ws.Range("Q6").Value = 0.07977417912033334
$ws.Range("R6").Value = 0.717967612083
$ws.Range("S6").Value = 0.0003119896219485922
$ws.Range("T6").Value = 0.0003119896219485922
$ws.Range("I7").Value = 0.2833335737960661
$ws.Range("J7").Value = 0.2833335737960661
$ws.Range("O7").Value = 0.00657695954769643
$ws.Range("P7").Value = 0.006576959547696431
$ws.Range("S7").Value = 0.001863473453360988
$ws.Range("T7").Value = 0.001863473453360988
$ws.Range("I8").Value = 0.2833335737960661
$ws.Range("J8").Value = 0.2833335737960661
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.008175
$ws.Range("N8").Value = 0.024525
$ws.Range("O8").Value = 0.0001396069650378791
$ws.Range("P8").Value = 0.0001396069650378791
$ws.Range("Q8").Value = 0.010114101825
$ws.Range("R8").Value = 0.09102691642500001
$ws.Range("S8").Value = 0.00003955534033100473
$ws.Range("T8").Value = 0.00003955534033100473
$ws.Range("I9").Value = 0.2833335737960661
$ws.Range("J9").Value = 0.2833335737960661
$ws.Range("M9").Value = 58.099467
$ws.Range("N9").Value = 174.298401
$ws.Range("O9").Value = 0.992182294579622
$ws.Range("P9").Value = 0.992182294579622
$ws.Range("Q9").Value = 71.88060247293301
$ws.Range("R9").Value = 646.925422256397
$ws.Range("S9").Value = 0.2811185553804255
$ws.Range("T9").Value = 0.2811185553804255
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1530633333333333
$ws.Range("H10").Value = 0.45919
$ws.Range("I10").Value = 0.03505335944376924
$ws.Range("J10").Value = 0.03505335944376924
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.06447966666666667
$ws.Range("N10").Value = 0.193439
$ws.Range("O10").Value = 0.001101138907643723
$ws.Range("P10").Value = 0.001101138907643722
$ws.Range("Q10").Value = 0.009869472712222223
$ws.Range("R10").Value = 0.08882525441
$ws.Range("S10").Value = 0.00003859861792715482
$ws.Range("T10").Value = 0.00003859861792715482
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1530633333333333
$ws.Range("H11").Value = 0.45919
$ws.Range("I11").Value = 0.03505335944376924
$ws.Range("J11").Value = 0.03505335944376924
$ws.Range("O11").Value = 0.00657695954769643
$ws.Range("P11").Value = 0.006576959547696431
$ws.Range("Q11").Value = 0.05894907748222222
$ws.Range("R11").Value = 0.53054169734
$ws.Range("S11").Value = 0.0002305445270725329
$ws.Range("T11").Value = 0.0002305445270725329
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1530633333333333
$ws.Range("H12").Value = 0.45919
$ws.Range("I12").Value = 0.03505335944376924
$ws.Range("J12").Value = 0.03505335944376924
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008175
$ws.Range("N12").Value = 0.024525
$ws.Range("O12").Value = 0.0001396069650378791
$ws.Range("P12").Value = 0.0001396069650378791
$ws.Range("Q12").Value = 0.00125129275
$ws.Range("R12").Value = 0.01126163475
$ws.Range("S12").Value = 0.0000048936931263265
$ws.Range("T12").Value = 0.0000048936931263265
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1530633333333333
$ws.Range("H13").Value = 0.45919
$ws.Range("I13").Value = 0.03505335944376924
$ws.Range("J13").Value = 0.03505335944376924
$ws.Range("M13").Value = 58.099467
$ws.Range("N13").Value = 174.298401
$ws.Range("O13").Value = 0.992182294579622
$ws.Range("P13").Value = 0.992182294579622
$ws.Range("Q13").Value = 8.89289808391
$ws.Range("R13").Value = 80.03608275519001
$ws.Range("S13").Value = 0.03477932260564322
$ws.Range("T13").Value = 0.03477932260564322
$ws.Range("G14").Value = 1.146373333333333
$ws.Range("H14").Value = 3.43912
$ws.Range("I14").Value = 0.2625333947391181
$ws.Range("J14").Value = 0.2625333947391181
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06447966666666667
$ws.Range("N14").Value = 0.193439
$ws.Range("O14").Value = 0.001101138907643723
$ws.Range("P14").Value = 0.001101138907643722
$ws.Range("Q14").Value = 0.07391777040888889
$ws.Range("R14").Value = 0.66525993368
$ws.Range("S14").Value = 0.0002890857355030307
$ws.Range("T14").Value = 0.0002890857355030307
$ws.Range("G15").Value = 1.146373333333333
$ws.Range("H15").Value = 3.43912
$ws.Range("I15").Value = 0.2625333947391181
$ws.Range("J15").Value = 0.2625333947391181
$ws.Range("O15").Value = 0.00657695954769643
$ws.Range("P15").Value = 0.006576959547696431
$ws.Range("Q15").Value = 0.4415012333688889
$ws.Range("R15").Value = 3.97351110032
$ws.Range("S15").Value = 0.001726671517118599
$ws.Range("T15").Value = 0.001726671517118599
$ws.Range("G16").Value = 1.146373333333333
$ws.Range("H16").Value = 3.43912
$ws.Range("I16").Value = 0.2625333947391181
$ws.Range("J16").Value = 0.2625333947391181
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.008175
$ws.Range("N16").Value = 0.024525
$ws.Range("O16").Value = 0.0001396069650378791
$ws.Range("P16").Value = 0.0001396069650378791
$ws.Range("Q16").Value = 0.009371602
$ws.Range("R16").Value = 0.084344418
$ws.Range("S16").Value = 0.00003665149046061978
$ws.Range("T16").Value = 0.00003665149046061978
$ws.Range("G17").Value = 1.146373333333333
$ws.Range("H17").Value = 3.43912
$ws.Range("I17").Value = 0.2625333947391181
$ws.Range("J17").Value = 0.2625333947391181
$ws.Range("M17").Value = 58.099467
$ws.Range("N17").Value = 174.298401
$ws.Range("O17").Value = 0.992182294579622
$ws.Range("P17").Value = 0.992182294579622
$ws.Range("Q17").Value = 66.60367964968
$ws.Range("R17").Value = 599.43311684712
$ws.Range("S17").Value = 0.2604809859960359
$ws.Range("T17").Value = 0.2604809859960359
